$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("エラーコード一覧")
$ws2 = $wb.Worksheets.Item("メッセージ一覧")

# --- Sheet1 (エラーコード一覧): add rows 10-12 ---
$ws1.Range("A10").Value = "0008"
$ws1.Range("B10").Value = "lỗi liên quan đến sản phẩm"

$ws1.Range("A11").Value = "0009"
$ws1.Range("B11").Value = "định dạng hình ảnh không đúng"

$ws1.Range("A12").Value = "0010"
$ws1.Range("B12").Value = "lưu hình ảnh thất bại"

$ws1.Range("A10:B12").NumberFormat = "@"

# --- Sheet2 (メッセージ一覧): fill rows 10-11 ---
$ws2.Range("A10").Value = "0008"
$ws2.Range("B10").Value = "W"
$ws2.Range("D10").Value = "Định dạng hình ảnh không chính xác"

$ws2.Range("A11").Value = "0009"
$ws2.Range("B11").Value = "W"
$ws2.Range("D11").Value = "{0} thất bại"

$ws2.Range("A10:A11").NumberFormat = "@"

# Selections to match final state
$ws1.Range("A13").Select()
$ws2.Range("D11").Select()
